$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
$tr1.Delete()
$tr1.Text = "A slide"

$tr2 = $s.Shapes.Item(4).TextFrame.TextRange
$tr2.Delete()
$tr2.Text = "Followed by a picture"
